# Apply the "Add files via upload" edit to the budget/database workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Workbook-level protection marker (adds <workbookProtection/> to workbook.xml) ---
$null = $wb.Protect()

# --- Replace the stray "hello" cell + the old Sunday/Cake/Drink/Buy 7-11 block ---
# (rows 30:34 in the original) with the new Thursday/Cake/Drink/Buy 7-11 block
# (rows 28:31).
$ws.Range("B30:K34").ClearContents()

# Row 28 - new "Thursday" entry
$ws.Range("B28").Value = "Thursday"
$ws.Range("C28").Value = 9
$ws.Range("D28").Value = 3
$ws.Range("E28").Value = 3
$ws.Range("F28").Value = "'06/18/2020"
$ws.Range("G28").Value = "'$601.28"
$ws.Range("H28").Value = "Food"
$ws.Range("I28").Value = "'$150.00"

# Row 29 - Cake
$ws.Range("H29").Value = "Cake"
$ws.Range("I29").Value = "'$150.00"

# Row 30 - Drink
$ws.Range("H30").Value = "Drink"
$ws.Range("I30").Value = "'$301.28"

# Row 31 - Buy 7-11
$ws.Range("H31").Value = "Buy 7-11"

# --- Sheet view: drop the scrolled/zoomed-in view, reselect A1 ---
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$null = $ws.Range("A1").Select()

Write-Output "edit applied"
